$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 64137.57
$ws.Range("I6").Value = 385
$ws.Range("K6").Value = 1155
$ws.Range("M6").Value = -1043

$ws.Range("H17").Value = 816.6
$ws.Range("J17").Value = 816.6
$ws.Range("L17").Value = 2449.8
$ws.Range("N17").Value = -2785.8

$ws.Range("H31").Value = 400.5
$ws.Range("I31").Value = 400.5
$ws.Range("K31").Value = 1201.5
$ws.Range("M31").Value = -971.5

$ws.Range("H40").Value = 991.6667
$ws.Range("I40").Value = 950
$ws.Range("K40").Value = 950
$ws.Range("M40").Value = -775

$ws.Range("H103").Value = 1175.2941
$ws.Range("I103").Value = 800
$ws.Range("J103").Value = 1711.4286
$ws.Range("K103").Value = 2400
$ws.Range("L103").Value = 5134.2858
$ws.Range("M103").Value = -1814
$ws.Range("N103").Value = -6306.2858

$ws.Range("H111").Value = 1435.25
$ws.Range("I111").Value = 1855.9286
$ws.Range("J111").Value = 846.3
$ws.Range("K111").Value = 5567.7858
$ws.Range("L111").Value = 2538.9
$ws.Range("M111").Value = -2500.7858
$ws.Range("N111").Value = -8672.9

$ws.Range("H112").Value = 1906.2
$ws.Range("J112").Value = 2116.9048
$ws.Range("L112").Value = 6350.714399999999
$ws.Range("N112").Value = -8566.714399999999

$ws.Range("H129").Value = 940.4474
$ws.Range("J129").Value = 1058.5
$ws.Range("L129").Value = 3175.5
$ws.Range("N129").Value = -13175.5

$ws.Range("H138").Value = 2489669
$ws.Range("I138").Value = 1695.9166
$ws.Range("J138").Value = 5378928
$ws.Range("K138").Value = 5087.7498
$ws.Range("L138").Value = 16136784
$ws.Range("M138").Value = 52.2502000000004
$ws.Range("N138").Value = -16147064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 3002
$ws.Range("I6").Value = 3002
$ws.Range("K6").Value = 3002
$ws.Range("M6").Value = -2829

$ws.Range("H32").Value = 2511.55
$ws.Range("I32").Value = 2047.1097
$ws.Range("J32").Value = 4627.3335
$ws.Range("K32").Value = 2047.1097
$ws.Range("L32").Value = 4627.3335
$ws.Range("M32").Value = -1760.1097
$ws.Range("N32").Value = -5201.3335

$ws.Range("H97").Value = 2084508.4
$ws.Range("I97").Value = 3126327
$ws.Range("J97").Value = 871.1
$ws.Range("K97").Value = 3126327
$ws.Range("L97").Value = 871.1
$ws.Range("M97").Value = -3125831
$ws.Range("N97").Value = -1863.1

$ws.Range("H110").Value = 400886.97
$ws.Range("I110").Value = 500679
$ws.Range("J110").Value = 1718.8
$ws.Range("K110").Value = 500679
$ws.Range("L110").Value = 1718.8
$ws.Range("M110").Value = -498634
$ws.Range("N110").Value = -5808.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 905.2222
$ws.Range("J94").Value = 1196.6666
$ws.Range("L94").Value = 1196.6666
$ws.Range("N94").Value = -2098.6666

$ws.Range("H99").Value = 971.9286
$ws.Range("J99").Value = 1059.6666
$ws.Range("L99").Value = 1059.6666
$ws.Range("N99").Value = -4055.6666

$ws.Range("H105").Value = 19233096
$ws.Range("I105").Value = 33335588
$ws.Range("K105").Value = 33335588
$ws.Range("M105").Value = -33333841

$ws.Range("H107").Value = 2938.2727
$ws.Range("I107").Value = 2702.3333
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 2702.3333
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -782.3332999999998
$ws.Range("N107").Value = -7840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -2684
$ws.Range("N32").Value = -3632

$ws.Range("H134").Value = 19893.31
$ws.Range("I134").Value = 1189.4131
$ws.Range("J134").Value = 91591.586
$ws.Range("K134").Value = 3568.2393
$ws.Range("L134").Value = 274774.758
$ws.Range("M134").Value = -1033.2393
$ws.Range("N134").Value = -279844.758

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 717
$ws.Range("I6").Value = 43.8
$ws.Range("J6").Value = 2400
$ws.Range("K6").Value = 131.4
$ws.Range("L6").Value = 7200
$ws.Range("M6").Value = -18.39999999999998
$ws.Range("N6").Value = -7426

$ws.Range("H7").Value = 357.07693
$ws.Range("I7").Value = 162.5
$ws.Range("K7").Value = 487.5
$ws.Range("M7").Value = -375.5

$ws.Range("H25").Value = 971.2857
$ws.Range("J25").Value = 966.5
$ws.Range("L25").Value = 2899.5
$ws.Range("N25").Value = -3237.5

$ws.Range("H30").Value = 971.2857
$ws.Range("J30").Value = 966.5
$ws.Range("L30").Value = 2899.5
$ws.Range("N30").Value = -3103.5

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H122").Value = 607.10254
$ws.Range("J122").Value = 1154
$ws.Range("L122").Value = 10386
$ws.Range("N122").Value = -15286

$ws.Range("H131").Value = 1001.2105
$ws.Range("I131").Value = 393
$ws.Range("J131").Value = 1093.3636
$ws.Range("K131").Value = 1179
$ws.Range("L131").Value = 3280.0908
$ws.Range("M131").Value = 3861
$ws.Range("N131").Value = -13360.0908

$ws.Range("H140").Value = 2314.6829
$ws.Range("I140").Value = 2329.1667
$ws.Range("J140").Value = 2303.348
$ws.Range("K140").Value = 6987.500100000001
$ws.Range("L140").Value = 6910.044
$ws.Range("M140").Value = -1807.500100000001
$ws.Range("N140").Value = -17270.044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2682.1
$ws.Range("I97").Value = 2601.25
$ws.Range("J97").Value = 3005.5
$ws.Range("K97").Value = 2601.25
$ws.Range("L97").Value = 3005.5
$ws.Range("M97").Value = -2105.25
$ws.Range("N97").Value = -3997.5

$ws.Range("H107").Value = 259.0909
$ws.Range("I107").Value = 130
$ws.Range("K107").Value = 130
$ws.Range("M107").Value = 1790

$ws.Range("H113").Value = 1205.238
$ws.Range("I113").Value = 925.8333
$ws.Range("J113").Value = 1577.7778
$ws.Range("K113").Value = 925.8333
$ws.Range("L113").Value = 1577.7778
$ws.Range("M113").Value = 1244.1667
$ws.Range("N113").Value = -5917.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 674.1458
$ws.Range("I46").Value = 610.2917
$ws.Range("J46").Value = 738
$ws.Range("K46").Value = 610.2917
$ws.Range("L46").Value = 738
$ws.Range("M46").Value = -422.2917
$ws.Range("N46").Value = -1114

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 45878.28
$ws.Range("I136").Value = 39757.73
$ws.Range("J136").Value = 53835
$ws.Range("K136").Value = 119273.19
$ws.Range("L136").Value = 161505
$ws.Range("M136").Value = -116723.19
$ws.Range("N136").Value = -166605
